$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "Sat Oct 12 23:34:50 EDT 2024",
    "Sat Oct 12 23:35:00 EDT 2024",
    "Sat Oct 12 23:35:11 EDT 2024",
    "Sat Oct 12 23:35:21 EDT 2024",
    "Sat Oct 12 23:35:31 EDT 2024",
    "Sat Oct 12 23:35:43 EDT 2024",
    "Sat Oct 12 23:35:53 EDT 2024",
    "Sat Oct 12 23:36:03 EDT 2024",
    "Sat Oct 12 23:36:13 EDT 2024",
    "Sat Oct 12 23:36:24 EDT 2024",
    "Sat Oct 12 23:36:34 EDT 2024",
    "Sat Oct 12 23:36:44 EDT 2024",
    "Sat Oct 12 23:36:54 EDT 2024",
    "Sat Oct 12 23:37:04 EDT 2024",
    "Sat Oct 12 23:37:14 EDT 2024",
    "Sat Oct 12 23:37:25 EDT 2024",
    "Sat Oct 12 23:37:36 EDT 2024",
    "Sat Oct 12 23:37:46 EDT 2024",
    "Sat Oct 12 23:37:56 EDT 2024",
    "Sat Oct 12 23:38:06 EDT 2024",
    "Sat Oct 12 23:38:16 EDT 2024",
    "Sat Oct 12 23:38:26 EDT 2024",
    "Sat Oct 12 23:38:38 EDT 2024",
    "Sat Oct 12 23:38:50 EDT 2024",
    "Sat Oct 12 23:39:00 EDT 2024",
    "Sat Oct 12 23:39:10 EDT 2024",
    "Sat Oct 12 23:39:20 EDT 2024",
    "Sat Oct 12 23:39:32 EDT 2024",
    "Sat Oct 12 23:39:42 EDT 2024",
    "Sat Oct 12 23:39:52 EDT 2024",
    "Sat Oct 12 23:40:02 EDT 2024",
    "Sat Oct 12 23:40:13 EDT 2024",
    "Sat Oct 12 23:40:23 EDT 2024",
    "Sat Oct 12 23:40:33 EDT 2024",
    "Sat Oct 12 23:40:44 EDT 2024",
    "Sat Oct 12 23:40:55 EDT 2024",
    "Sat Oct 12 23:41:05 EDT 2024",
    "Sat Oct 12 23:41:15 EDT 2024",
    "Sat Oct 12 23:41:25 EDT 2024",
    "Sat Oct 12 23:41:35 EDT 2024",
    "Sat Oct 12 23:41:45 EDT 2024",
    "Sat Oct 12 23:41:55 EDT 2024",
    "Sat Oct 12 23:42:06 EDT 2024",
    "Sat Oct 12 23:42:16 EDT 2024",
    "Sat Oct 12 23:42:27 EDT 2024",
    "Sat Oct 12 23:42:37 EDT 2024",
    "Sat Oct 12 23:42:47 EDT 2024"
)

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}

